$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap columns A, Q, R, AC between row 12 and row 13
$cols = "A", "Q", "R", "AC"
foreach ($col in $cols) {
    $cellTop = $ws.Range($col + "12")
    $cellBottom = $ws.Range($col + "13")
    $tmp = $cellTop.Value2
    $cellTop.Value2 = $cellBottom.Value2
    $cellBottom.Value2 = $tmp
}

# Swap columns A, Q, R, AC between row 15 and row 16
foreach ($col in $cols) {
    $cellTop = $ws.Range($col + "15")
    $cellBottom = $ws.Range($col + "16")
    $tmp = $cellTop.Value2
    $cellTop.Value2 = $cellBottom.Value2
    $cellBottom.Value2 = $tmp
}
